# Updates the cryptos price/volume table (columns B-E, rows 2-51) with the
# latest scrape values. Some Price (column D) values look like plain numbers
# (e.g. "1.002", "1.0000") which Excel would otherwise silently coerce to a
# Number on assignment; for those we force the cell to Text format first so
# the literal string (including trailing zeros) is preserved, matching the
# source feed's formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "24.427.42"
$ws.Range("E2").Value = "  -1.81%  "

# Row 3
$ws.Range("D3").Value = "1.650.45"
$ws.Range("E3").Value = "  -4.22%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -1.38%  "

# Row 5
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.33"
$ws.Range("E5").Value = "  -2.13%  "

# Row 6
$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9988"
$ws.Range("E6").Value = "  -0.93%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3627"
$ws.Range("E7").Value = "  -4.38%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.39"
$ws.Range("E8").Value = "  -4.68%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3281"
$ws.Range("E9").Value = "  -6.82%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.117"
$ws.Range("E10").Value = "  -6.76%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06911"
$ws.Range("E11").Value = "  -8.08%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9997"
$ws.Range("E12").Value = "  -0.93%  "

# Row 13
$ws.Range("E13").Value = "  -7.16%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.12"
$ws.Range("E14").Value = "  -8.84%  "

# Row 15
$ws.Range("D15").Value = "1.647.16"
$ws.Range("E15").Value = "  -5.25%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.552"
$ws.Range("E16").Value = "  -6.84%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001045"
$ws.Range("E17").Value = "  -7.54%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06489"
$ws.Range("E18").Value = "  -3.29%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.0000"
$ws.Range("E19").Value = "  -0.88%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "76.95"
$ws.Range("E20").Value = "  -9.51%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.907"
$ws.Range("E21").Value = "  -8.12%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.72"
$ws.Range("E22").Value = "  -9.77%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.18"
$ws.Range("E23").Value = "  -7.90%  "

# Row 24
$ws.Range("D24").Value = "24.393.82"
$ws.Range("E24").Value = "  -2.20%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.421"
$ws.Range("E25").Value = "  -1.17%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.330"
$ws.Range("E26").Value = "  -17.29%  "

# Row 27
$ws.Range("E27").Value = "  -4.25%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.46"
$ws.Range("E28").Value = "  -10.42%  "

# Row 29
$ws.Range("D29").Value = "1.829.73"
$ws.Range("E29").Value = "  -4.97%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.52"
$ws.Range("E30").Value = "  -6.02%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.153"
$ws.Range("E31").Value = "  -2.50%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.044"
$ws.Range("E32").Value = "  -4.85%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.578"
$ws.Range("E33").Value = "  -19.56%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08320"
$ws.Range("E34").Value = "  -4.83%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.667"
$ws.Range("E35").Value = "  -7.49%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.26"
$ws.Range("E36").Value = "  -11.36%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.205"
$ws.Range("E37").Value = "  -7.63%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06030"
$ws.Range("E38").Value = "  -8.65%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02209"
$ws.Range("E39").Value = "  -10.18%  "

# Row 40
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.256"
$ws.Range("E40").Value = "  -9.87%  "

# Row 41
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.202"
$ws.Range("E41").Value = "  -5.59%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2036"
$ws.Range("E42").Value = "  -8.53%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9991"
$ws.Range("E43").Value = "  -0.94%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5839"
$ws.Range("E44").Value = "  -10.13%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.718"
$ws.Range("E45").Value = "  -3.85%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.69"
$ws.Range("E46").Value = "  -9.25%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5566"
$ws.Range("E47").Value = "  -10.34%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "121.53"
$ws.Range("E48").Value = "  -6.22%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.931"
$ws.Range("E49").Value = "  -10.80%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06892"
$ws.Range("E50").Value = "  -5.59%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.84"
$ws.Range("E51").Value = "  -8.01%  "
